$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Type of Model" column label (A1); B1 "Score" already present
$ws.Range("A1").Value = "Type of Model"

# Existing Linear Regression row: update score value
$ws.Range("B2").Value = 0.23

# New rows for additional models
$ws.Range("A3").Value = "Logistic Regression"
$ws.Range("B3").Value = 0.36
$ws.Range("B3").NumberFormat = "0%"

$ws.Range("C3").Value = "82% if in order"
$ws.Range("C3").NumberFormat = "0%"

$ws.Range("A4").Value = "Gradient Descent"
$ws.Range("A5").Value = "KNN"
$ws.Range("A6").Value = "SVM"
$ws.Range("A7").Value = "Trees"

# Column A width (characters) - tuned so serialized width rounds to 17
$ws.Range("A1").ColumnWidth = 16.1667

# Selection moves to C4 as in the edited workbook
[void]$ws.Range("C4").Select()
